$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Frontend Bugs")

# Row 31
$ws.Range('A31').Value = 'A-030'
$ws.Range('B31').Value = 'Admin Landing Page'
$ws.Range('C31').Value = 'Chrome'
$ws.Range('D31').Value = 'Responsiveness'
$ws.Range('E31').Value = 'The Active Users and Response Time cards on Admin Setting image must not be visible on small screens'
$ws.Range('F31').Value = 'either remove them or make some change for maing it responsive'
$ws.Range('G31').Value = 'those 2 cards must be static beneath the setting image instead of floating on it'
$ws.Range('H31').Value = 'those 2 cards are floating on the setting image'
$ws.Range('I31').Value = 'medium'
$ws.Range('J31').Value = 'to do'
$ws.Hyperlinks.Add($ws.Range('K31'), 'bugs\screenshots\Landing-Page-Admin.png')
$ws.Range('K31').Style = "Hyperlink"
$ws.Range('L31').Value = '2nd Feb, 2026'
$ws.Range('M31').Value = '5th Feb, 2026'
$ws.Range('O31').Value = 'Muhammad Noman'

# Row 32
$ws.Range('A32').Value = 'A-031'
$ws.Range('B32').Value = 'Admin Landing Page'
$ws.Range('C32').Value = 'Chrome'
$ws.Range('D32').Value = 'Exagerated Lies'
$ws.Range('E32').Value = 'Stop presenting lies'
$ws.Range('F32').Value = 'remove ip whitelisting, automatic logout'
$ws.Range('G32').Value = 'correct info'
$ws.Range('H32').Value = 'incorrect info'
$ws.Range('I32').Value = 'high'
$ws.Range('J32').Value = 'to do'
$ws.Hyperlinks.Add($ws.Range('K32'), 'bugs\screenshots\landing-admin-details.png')
$ws.Range('K32').Style = "Hyperlink"
$ws.Range('L32').Value = '2nd Feb, 2026'
$ws.Range('M32').Value = '5th Feb, 2026'
$ws.Range('O32').Value = 'Muhammad Noman'

# Row 33 (B filled in later, after row 34)
$ws.Range('A33').Value = 'C-001'
$ws.Range('C33').Value = 'Chrome'
$ws.Range('D33').Value = 'UX'
$ws.Range('E33').Value = 'Though we have a single landing page for both customers and doctors it does not mean we are free to place any button anywhere. See how bad UX will be when a customer will see this login as doctor button. Just show a single login button on the top of the page at th very right top corner. and remove get started button from right top corner. further in hero section just show a get started button. remove login as doctor button.'
$ws.Range('F33').Value = 'update this landing page'
$ws.Range('G33').Value = 'as mentioned in description'
$ws.Range('H33').Value = 'current wrong state'
$ws.Range('I33').Value = 'medium'
$ws.Range('J33').Value = 'to do'
$ws.Hyperlinks.Add($ws.Range('K33'), 'bugs\screenshots\customer-doctor-landing.png')
$ws.Range('K33').Style = "Hyperlink"
$ws.Range('L33').Value = '2nd Feb, 2026'
$ws.Range('M33').Value = '5th Feb, 2026'
$ws.Range('O33').Value = 'Muhammad Noman'

# Row 34
$ws.Range('A34').Value = 'C-002'
$ws.Range('B34').Value = 'Customer Doctor Landing Page'
$ws.Range('C34').Value = 'Chrome'
$ws.Range('D34').Value = 'UX'
$ws.Range('E34').Value = 'Remove buttons in for doctors dection and make the heading for doctors more prominent'
$ws.Range('F34').Value = 'same as of description'
$ws.Range('G34').Value = 'just details'
$ws.Range('H34').Value = 'buttons are also shown'
$ws.Range('I34').Value = 'medium'
$ws.Range('J34').Value = 'to do'
$ws.Hyperlinks.Add($ws.Range('K34'), 'bugs\screenshots\doctors-landing.png')
$ws.Range('K34').Style = "Hyperlink"
$ws.Range('L34').Value = '2nd Feb, 2026'
$ws.Range('M34').Value = '5th Feb, 2026'
$ws.Range('O34').Value = 'Muhammad Noman'

# Back-fill B33 (reuses "Customer Doctor Landing Page" created above)
$ws.Range('B33').Value = 'Customer Doctor Landing Page'

# Row 35
$ws.Range('A35').Value = 'C-003'
$ws.Range('B35').Value = 'Customer Doctor Landing Page'
$ws.Range('C35').Value = 'Chrome'
$ws.Range('D35').Value = 'UX'
$ws.Range('E35').Value = 'Wrong date in footer'
$ws.Range('F35').Value = 'instead of hardcoding the date use Date object and get year from that object'
$ws.Range('G35').Value = 'current year'
$ws.Range('H35').Value = 'hard coded year'
$ws.Range('I35').Value = 'low'
$ws.Range('J35').Value = 'to do'
$ws.Hyperlinks.Add($ws.Range('K35'), 'bugs\screenshots\landing date.png')
$ws.Range('K35').Style = "Hyperlink"
$ws.Range('L35').Value = '2nd Feb, 2026'
$ws.Range('M35').Value = '5th Feb, 2026'
$ws.Range('O35').Value = 'Muhammad Noman'

# Row 36
$ws.Range('A36').Value = 'C-004'
$ws.Range('B36').Value = 'Customer Doctor Landing Page'
$ws.Range('C36').Value = 'Chrome'
$ws.Range('D36').Value = 'Responsiveness'
$ws.Range('E36').Value = 'philbox icons size is too large'
$ws.Range('F36').Value = 'update tailwindd classes'
$ws.Range('G36').Value = 'it must be responsive'
$ws.Range('H36').Value = 'irresponsive'
$ws.Range('I36').Value = 'medium'
$ws.Range('J36').Value = 'to do'
$ws.Hyperlinks.Add($ws.Range('K36'), 'bugs\screenshots\landing-doctor-foter-icon.png')
$ws.Range('K36').Style = "Hyperlink"
$ws.Range('L36').Value = '2nd Feb, 2026'
$ws.Range('M36').Value = '5th Feb, 2026'
$ws.Range('O36').Value = 'Muhammad Noman'

# Final selection matches the end-state sheet view (cursor parked on the next empty row)
$ws.Range('A37').Select()

